$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for new column H
$ws.Range("H1").Value = "Added_cf_num"
$ws.Range("H1").Font.Bold = $true

# Formula column: H3 standalone, H4:H35 filled via autofill (becomes shared formula)
$ws.Range("H3").Formula = "=G3-27155"
$ws.Range("H4:H35").Formula = "=G4-27155"

# Autofit column H width to match bestFit behavior
$ws.Columns("H").AutoFit()

# Update selection to H3 as last active cell
$ws.Range("H3").Select()
